# Comentario Documento virtuaDlab virtual 3
#
# Splits the single paragraph "Esto es una prueba de git" into a run
# that keeps "Esto es una prueba de " plus a spell-check-flagged run for
# "git", then appends a blank paragraph and a new paragraph
# "Prueba de gist 2" (with "gist" flagged by the spell checker), moving
# the _GoBack bookmark to the very end of the new content.

$d = $word.ActiveDocument

# The whole (only) paragraph in the document - replace its contents
# (everything after <w:pPr>) with the new run/proofErr/bookmark layout,
# and append the two new paragraphs that follow it. The paragraph's own
# <w:p ...rsid...> and <w:pPr> are left completely untouched because we
# only touch the Range that covers the paragraph's text content.
$range = $d.Paragraphs(1).Range

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="005C17F9" w:rsidRPr="00C42543" w:rsidRDefault="00C42543"><w:pPr><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">Esto es una prueba de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>git</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="es-419"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve">Prueba de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t>gist</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="es-419"/></w:rPr><w:t xml:space="preserve"> 2</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$range.InsertXML($xml)
